# end_of_day_reports_input.xlsx edit
# - "input" sheet (DataHub report rows 506-534): fix column ordering / add
#   two new derived-date rows (EpisodeDate_imputed, OnsetDate_imputed) and
#   flip the PT row's clean_it flag to "no".
# - Active tab moves from "lbls" back to "input"; selections updated on
#   both "input" and "lbls" sheets.

$wb = $excel.ActiveWorkbook
$wsInput = $wb.Worksheets.Item("input")
$wsLbls  = $wb.Worksheets.Item("lbls")

# --- 1. Fix the "PT" row's clean_it column (row 506, col C: yes -> no) ---
$wsInput.Cells.Item(506, 3).Value = "no"

# --- 2. EpisodeDate (row 507): set case_standard (col F) to str_to_lower ---
$wsInput.Cells.Item(507, 6).Value = "str_to_lower"

# --- 3. Insert a new row after EpisodeDate (507) for EpisodeDate_imputed,
#        cloning EpisodeDate's row values but with the new column name ---
$wsInput.Rows.Item(508).Insert()
$wsInput.Cells.Item(508, 1).Value = "DataHub"
$wsInput.Cells.Item(508, 2).Value = "EpisodeDate_imputed"
$wsInput.Cells.Item(508, 3).Value = "no"
$wsInput.Cells.Item(508, 4).Value = "-"
$wsInput.Cells.Item(508, 5).Value = " "
$wsInput.Cells.Item(508, 6).Value = "str_to_lower"

# --- 4. OnsetDate2 (previously row 516, now shifted to row 517 after the
#        insert above): rename to OnsetDate and set case_standard ---
$wsInput.Cells.Item(517, 2).Value = "OnsetDate"
$wsInput.Cells.Item(517, 6).Value = "str_to_lower"

# --- 5. Insert a new row after OnsetDate (517) for OnsetDate_imputed ---
$wsInput.Rows.Item(518).Insert()
$wsInput.Cells.Item(518, 1).Value = "DataHub"
$wsInput.Cells.Item(518, 2).Value = "OnsetDate_imputed"
$wsInput.Cells.Item(518, 3).Value = "no"
$wsInput.Cells.Item(518, 4).Value = "-"
$wsInput.Cells.Item(518, 5).Value = " "
$wsInput.Cells.Item(518, 6).Value = "str_to_lower"

# --- 5b. RecoveryDate2 (now shifted down to row 535): case_standard -> str_to_lower ---
$wsInput.Cells.Item(535, 6).Value = "str_to_lower"

# --- 6. Re-apply the autofilter over the new, larger range and update the
#        hidden _FilterDatabase defined name to match ---
$wsInput.AutoFilterMode = $false
$wsInput.Range("A1:F536").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name() -eq "input!_FilterDatabase") {
        $n.RefersTo = "=input!`$A`$1:`$F`$536"
    }
}

# --- 7. Update view state: "lbls" loses the active-tab flag and gets a new
#        selection; "input" becomes the active sheet with the new selection ---
$wsLbls.Range("H7:H8").Select()
$wsInput.Activate()
$wsInput.Range("B507").Select()
